$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = 4201.730337248195
$ws.Range("C3").Value  = 4201.730337248195
$ws.Range("C4").Value  = 4153.802161834338
$ws.Range("C5").Value  = 4098.736683771019
$ws.Range("C6").Value  = 4098.736683771019
$ws.Range("C7").Value  = 4098.736683771019
$ws.Range("C8").Value  = 4098.736683771019
$ws.Range("C9").Value  = 3963.054608463689
$ws.Range("C10").Value = 3963.054608463689
$ws.Range("C11").Value = 3963.054608463689
$ws.Range("C12").Value = 3963.054608463689
